$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.860.48'
$ws.Range("E2").Value = '  +2.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.090.13'
$ws.Range("E3").Value = '  +5.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.80'
$ws.Range("E5").Value = '  +2.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.75'
$ws.Range("E6").Value = '  +6.29%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.085.31'
$ws.Range("E8").Value = '  +5.25%  '

$ws.Range("E9").Value = '  +1.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.63'
$ws.Range("E10").Value = '  -1.53%  '

$ws.Range("E11").Value = '  +3.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.481'
$ws.Range("E12").Value = '  +5.80%  '

$ws.Range("E13").Value = '  +1.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.37'
$ws.Range("E14").Value = '  +6.33%  '

$ws.Range("E15").Value = '  -0.50%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.598.10'
$ws.Range("E16").Value = '  +5.12%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.865.00'
$ws.Range("E17").Value = '  +2.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.20'
$ws.Range("E18").Value = '  +4.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.090.00'
$ws.Range("E19").Value = '  +5.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.18'
$ws.Range("E20").Value = '  +7.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '466.34'
$ws.Range("E21").Value = '  +4.96%  '

$ws.Range("E22").Value = '  +4.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("E23").Value = '  +4.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.77'
$ws.Range("E24").Value = '  +2.12%  '

$ws.Range("E25").Value = '  +6.69%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.07'
$ws.Range("E26").Value = '  +8.02%  '

$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.00'
$ws.Range("E29").Value = '  -0.88%  '

$ws.Range("E30").Value = '  -0.44%  '

$ws.Range("E31").Value = '  +4.08%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0000102'
$ws.Range("E32").Value = '  +1.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.20'
$ws.Range("E33").Value = '  +4.29%  '

$ws.Range("E34").Value = '  +3.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  +3.70%  '

$ws.Range("E37").Value = '  +3.33%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.11'
$ws.Range("E38").Value = '  +6.80%  '

$ws.Range("B39").Value = 'Arweave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '46.97'
$ws.Range("E39").Value = '  +4.82%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.319'
$ws.Range("E40").Value = '  +6.91%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.29'
$ws.Range("E41").Value = '  +1.15%  '

$ws.Range("E42").Value = '  +1.52%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.68'
$ws.Range("E43").Value = '  +2.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.83'
$ws.Range("E44").Value = '  +0.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0361'
$ws.Range("E45").Value = '  +2.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '383.37'
$ws.Range("E46").Value = '  -0.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.778.93'
$ws.Range("E47").Value = '  +2.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.15'
$ws.Range("E48").Value = '  +1.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.88'
$ws.Range("E50").Value = '  +6.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.23'
$ws.Range("E51").Value = '  +2.32%  '
